# Uruguay Primera División - base update (29-02-2024 07:50)
# 1) Rows 117-120: odds feed re-sync (data shifted a row; id (col A) kept fixed)
# 2) Rows 136-142: 7 new upcoming fixtures appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update existing rows 117-120 (all columns except A, the row id, change)
# ---------------------------------------------------------------------------
$updated = @{
    117 = @{ B=7013409; F="Nacional De Football"; G="Torque"; H=1; I=1; J="D";
             K=1.666; L=3.9;  M=4.5;  N=1.615; O=4;     P=4.75;
             Q=-0.75; R=1.8;   S=2.05;  T=2.75; U=1.95;  V=1.9;
             W=-1;    X=3;     Y=-1;    Z=-1;   AA=1.05; AB=-1;   AC=0.8999999999999999 }
    118 = @{ B=7013885; F="La Luz"; G="Atletico Fenix Montevideo"; H=0; I=2; J="A";
             K=3;     L=3;     M=2.4;   N=2.9;  O=2.75;  P=2.6;
             Q=0;     R=2.025; S=1.825; T=2;    U=2.025; V=1.825;
             W=-1;    X=-1;    Y=1.6;   Z=-1;   AA=0.825;AB=0;    AC=-0 }
    119 = @{ B=7013702; F="Defensor Sporting"; G="Danubio"; H=0; I=2; J="A";
             K=1.8;   L=3.6;   M=4.2;   N=1.8;  O=3.6;   P=4.2;
             Q=-0.75; R=2.05;  S=1.8;   T=2.25; U=1.85;  V=2;
             W=-1;    X=-1;    Y=3.2;   Z=-1;   AA=0.8;  AB=-0.5; AC=0.5 }
    120 = @{ B=7013886; F="Racing Club de Montevideo"; G="Cerro"; H=0; I=1; J="A";
             K=2.25;  L=3.1;   M=3.25;  N=2.25; O=2.875; P=3.5;
             Q=-0.25; R=1.95;  S=1.9;   T=2;    U=1.925; V=1.925;
             W=-1;    X=-1;    Y=2.5;   Z=-1;   AA=0.8999999999999999; AB=-1; AC=0.925 }
}

$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

foreach ($r in $updated.Keys) {
    $vals = $updated[$r]
    foreach ($c in $cols) {
        $ws.Range($c + $r).Value = $vals[$c]
    }
}

# ---------------------------------------------------------------------------
# 2) Append 7 new fixtures as rows 136-142
# ---------------------------------------------------------------------------
$newRows = @{
    136 = @{ A=134; B=7870599; C="Uruguay Primera División"; D="Uruguay Apertura"; E=45352.70833333334;  F="Cerro";                      G="Club Atletico Progreso";
             K=2.05;  L=3.1; M=3.75; N=2;     O=3.1; P=4;
             Q=-0.5;  R=2.05;  S=1.8;   T=2.25; U=2.025; V=1.825;
             W=0; X=0; Y=0; Z=0; AA=0 }
    137 = @{ A=135; B=7870600; C="Uruguay Primera División"; D="Uruguay Apertura"; E=45352.8125;         F="Miramar Misiones";          G="Cerro Largo";
             K=2.875; L=3;   M=2.5;  N=2.9;   O=3;   P=2.45;
             Q=0;     R=2.125; S=1.75;  T=2;    U=1.85;  V=2;
             W=0; X=0; Y=0; Z=0; AA=0 }
    138 = @{ A=136; B=7870604; C="Uruguay Primera División"; D="Uruguay Apertura"; E=45353.41666666666;  F="Racing Club de Montevideo";  G="Defensor Sporting";
             K=3.75;  L=3.2; M=2;    N=3.8;   O=3.2; P=2;
             Q=0.5;   R=1.8;   S=2.05;  T=2.25; U=2;     V=1.85;
             W=0; X=0; Y=0; Z=0; AA=0 }
    139 = @{ A=137; B=7870605; C="Uruguay Primera División"; D="Uruguay Apertura"; E=45353.70833333334;  F="Danubio";                    G="CA River Plate";
             K=2.1;   L=3.1; M=3.6;  N=2.1;   O=3.1; P=3.6;
             Q=-0.25; R=1.825; S=2.025; T=2;    U=1.9;   V=1.95;
             W=0; X=0; Y=0; Z=0; AA=0 }
    140 = @{ A=138; B=7870606; C="Uruguay Primera División"; D="Uruguay Apertura"; E=45353.83333333334;  F="Nacional De Football";       G="Montevideo Wanderers";
             K=1.533; L=4;   M=5.5;  N=1.571; O=4;   P=5.25;
             Q=-1;    R=2.025; S=1.825; T=2.5;  U=1.95;  V=1.9;
             W=0; X=0; Y=0; Z=0; AA=0 }
    141 = @{ A=139; B=7870601; C="Uruguay Primera División"; D="Uruguay Apertura"; E=45354.41666666666;  F="Atletico Fenix Montevideo";  G="Boston River";
             K=2.375; L=3.1; M=3;    N=2.3;   O=3.1; P=3.1;
             Q=-0.25; R=2.025; S=1.825; T=2.25; U=2.025; V=1.825;
             W=0; X=0; Y=0; Z=0; AA=0 }
    142 = @{ A=140; B=7870602; C="Uruguay Primera División"; D="Uruguay Apertura"; E=45354.70833333334;  F="Liverpool Montevideo";       G="Deportivo Maldonado";
             K=1.909; L=3.3; M=4;    N=1.833; O=3.4; P=4.2;
             Q=-0.5;  R=1.875; S=1.975; T=2.5;  U=2.025; V=1.825;
             W=0; X=0; Y=0; Z=0; AA=0 }
}

$newCols = @("A","B","C","D","E","F","G","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

# Sort ascending so rows are created in order
$sortedRowNums = $newRows.Keys | Sort-Object

foreach ($r in $sortedRowNums) {
    # Carry the id (col A) and date (col E) number-format/border styling down from
    # the existing data block (row 135) before writing values into the new row.
    $ws.Range("A135").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("E135").Copy()
    $ws.Range("E$r").PasteSpecial(-4122)

    $vals = $newRows[$r]
    foreach ($c in $newCols) {
        $ws.Range($c + $r).Value = $vals[$c]
    }
}
